# Refatorando o consolidador para modelo ETL
# Update the absenteeism data rows (2-11) with the new ETL-sourced values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 50835
$ws.Range("B2").Value = "Sra. Ana Lívia Rezende"
$ws.Range("D2").Value = "Problemas pessoais"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 45087
$ws.Range("G2").Value = 8658.870000000001

# Row 3
$ws.Range("A3").Value = 52022
$ws.Range("B3").Value = "Marina Nogueira"
$ws.Range("C3").Value = "TI"
$ws.Range("D3").Value = "Problemas pessoais"
$ws.Range("E3").Value = 5
$ws.Range("G3").Value = 10599.4

# Row 4
$ws.Range("A4").Value = 34568
$ws.Range("B4").Value = "Srta. Lavínia Mendes"
$ws.Range("C4").Value = "Engenharia"
$ws.Range("D4").Value = "Consulta médica"
$ws.Range("E4").Value = 7
$ws.Range("F4").Value = 45094
$ws.Range("G4").Value = 2859.45

# Row 5
$ws.Range("A5").Value = 31453
$ws.Range("B5").Value = "Lorenzo Costa"
$ws.Range("C5").Value = "Atendimento ao Cliente"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 45091
$ws.Range("G5").Value = 5193.61

# Row 6
$ws.Range("A6").Value = 4410
$ws.Range("B6").Value = "Juliana Barbosa"
$ws.Range("C6").Value = "P&D"
$ws.Range("D6").Value = "Consulta médica"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 45088
$ws.Range("G6").Value = 12273.55

# Row 7
$ws.Range("A7").Value = 36679
$ws.Range("B7").Value = "Srta. Caroline Duarte"
$ws.Range("C7").Value = "Jurídico"
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = 45089
$ws.Range("G7").Value = 9222.129999999999

# Row 8
$ws.Range("A8").Value = 84454
$ws.Range("B8").Value = "Dra. Letícia Moura"
$ws.Range("C8").Value = "Recursos Humanos"
$ws.Range("D8").Value = "Doença"
$ws.Range("E8").Value = 6
$ws.Range("F8").Value = 45085
$ws.Range("G8").Value = 11769.08

# Row 9
$ws.Range("A9").Value = 1507
$ws.Range("B9").Value = "Fernanda Fernandes"
$ws.Range("C9").Value = "Marketing"
$ws.Range("D9").Value = "Consulta médica"
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 8199.27

# Row 10
$ws.Range("A10").Value = 58115
$ws.Range("B10").Value = "Srta. Isadora Fernandes"
$ws.Range("D10").Value = "Viagem de negócios"
$ws.Range("E10").Value = 8
$ws.Range("F10").Value = 45089
$ws.Range("G10").Value = 9262.24

# Row 11
$ws.Range("A11").Value = 98595
$ws.Range("B11").Value = "Sr. Enrico Souza"
$ws.Range("C11").Value = "Financeiro"
$ws.Range("D11").Value = "Doença"
$ws.Range("E11").Value = 7
$ws.Range("F11").Value = 45083
$ws.Range("G11").Value = 3451.41
